$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at X so the old m_elu/m_feu data shifts: old W(m_elu) stays,
# old X (m_feu) moves to Y, and the newly inserted X is blank for the new m_elu_2 formula column.
$ws.Columns("X").Insert()

# Split header "m_elu" into "m_elu_1" (W) and "m_elu_2" (X); Y keeps "m_feu" from the shift.
$ws.Range("W1").Value = "m_elu_1"
$ws.Range("X1").Value = "m_elu_2"

# Recompute m_els_1 / m_els_2 (U/V) for each data row, and fill new m_elu_2 (X) = U + V.
$ws.Range("U2").Value = 29.253119999999999
$ws.Range("V2").Value = 8.8857599999999994
$ws.Range("X2").Formula = "=U2+V2"
$ws.Range("U3").Value = 80.574999999999989
$ws.Range("V3").Value = 24.474999999999998
$ws.Range("X3").Formula = "=U3+V3"
$ws.Range("U4").Value = 30.75328
$ws.Range("V4").Value = 9.3414400000000004
$ws.Range("X4").Formula = "=U4+V4"
$ws.Range("U5").Value = 58.014000000000003
$ws.Range("V5").Value = 17.622000000000003
$ws.Range("X5").Formula = "=U5+V5"
$ws.Range("U6").Value = 30.0032
$ws.Range("V6").Value = 9.1135999999999999
$ws.Range("X6").Formula = "=U6+V6"
$ws.Range("U7").Value = 96.69
$ws.Range("V7").Value = 29.369999999999997
$ws.Range("X7").Formula = "=U7+V7"
$ws.Range("U8").Value = 52.505600000000001
$ws.Range("V8").Value = 15.9488
$ws.Range("X8").Formula = "=U8+V8"
$ws.Range("U9").Value = 70.906000000000006
$ws.Range("V9").Value = 21.538
$ws.Range("X9").Formula = "=U9+V9"
$ws.Range("U10").Value = 5.8506239999999998
$ws.Range("V10").Value = 1.7771519999999998
$ws.Range("X10").Formula = "=U10+V10"
$ws.Range("U11").Value = 19.337999999999997
$ws.Range("V11").Value = 5.8739999999999997
$ws.Range("X11").Formula = "=U11+V11"
$ws.Range("U12").Value = 6.7507199999999994
$ws.Range("V12").Value = 2.0505599999999999
$ws.Range("X12").Formula = "=U12+V12"
$ws.Range("U13").Value = 70.906000000000006
$ws.Range("V13").Value = 21.538
$ws.Range("X13").Formula = "=U13+V13"
$ws.Range("U14").Value = 42.754560000000005
$ws.Range("V14").Value = 12.986880000000001
$ws.Range("X14").Formula = "=U14+V14"
$ws.Range("U15").Value = 93.466999999999999
$ws.Range("V15").Value = 28.390999999999998
$ws.Range("X15").Formula = "=U15+V15"
$ws.Range("U16").Value = 33.753599999999999
$ws.Range("V16").Value = 10.252800000000001
$ws.Range("X16").Formula = "=U16+V16"
$ws.Range("U17").Value = 27.752959999999998
$ws.Range("V17").Value = 8.4300800000000002
$ws.Range("X17").Formula = "=U17+V17"
$ws.Range("U18").Value = 96.69
$ws.Range("V18").Value = 29.369999999999997
$ws.Range("X18").Formula = "=U18+V18"
$ws.Range("U19").Value = 56.256
$ws.Range("V19").Value = 17.088000000000001
$ws.Range("X19").Formula = "=U19+V19"
$ws.Range("U20").Value = 70.906000000000006
$ws.Range("V20").Value = 21.538
$ws.Range("X20").Formula = "=U20+V20"
$ws.Range("U21").Value = 2.1002239999999999
$ws.Range("V21").Value = 0.63795199999999985
$ws.Range("X21").Formula = "=U21+V21"
$ws.Range("U22").Value = 19.337999999999997
$ws.Range("V22").Value = 5.8739999999999997
$ws.Range("X22").Formula = "=U22+V22"
$ws.Range("U23").Value = 5.7006079999999999
$ws.Range("V23").Value = 1.7315839999999998
$ws.Range("X23").Formula = "=U23+V23"
$ws.Range("U24").Value = 70.906000000000006
$ws.Range("V24").Value = 21.538
$ws.Range("X24").Formula = "=U24+V24"
$ws.Range("U25").Value = 56.256
$ws.Range("V25").Value = 17.088000000000001
$ws.Range("X25").Formula = "=U25+V25"
$ws.Range("U26").Value = 93.466999999999999
$ws.Range("V26").Value = 28.390999999999998
$ws.Range("X26").Formula = "=U26+V26"
$ws.Range("U27").Value = 27.752959999999998
$ws.Range("V27").Value = 8.4300800000000002
$ws.Range("X27").Formula = "=U27+V27"

# Row 22 had its m_elu_1 (W) manually revised from 33 to 190.
$ws.Range("W22").Value = 190

$ws.Range("AB9").Select()
